$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "60.948.11"
$ws.Range("E2").Value = "  +1.07%  "

# Row 3
$ws.Range("D3").Value = "2.350.06"
$ws.Range("E3").Value = "  -1.55%  "

# Row 4
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.27%  "

# Row 5
$ws.Range("D5").Value = "545.12"
$ws.Range("E5").Value = "  +1.17%  "

# Row 6
$ws.Range("D6").Value = "136.76"
$ws.Range("E6").Value = "  -1.84%  "

# Row 7
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.13%  "

# Row 8
$ws.Range("D8").Value = "0.527"
$ws.Range("E8").Value = "  -8.14%  "

# Row 9
$ws.Range("D9").Value = "2.345.98"
$ws.Range("E9").Value = "  -1.67%  "

# Row 10
$ws.Range("D10").Value = "0.105"
$ws.Range("E10").Value = "  +1.04%  "

# Row 11
$ws.Range("E11").Value = "  +1.54%  "

# Row 12
$ws.Range("D12").Value = "5.31"
$ws.Range("E12").Value = "  +0.00%  "

# Row 13
$ws.Range("D13").Value = "0.340"
$ws.Range("E13").Value = "  +0.25%  "

# Row 14
$ws.Range("D14").Value = "24.62"
$ws.Range("E14").Value = "  -2.29%  "

# Row 15
$ws.Range("D15").Value = "2.768.72"
$ws.Range("E15").Value = "  -1.71%  "

# Row 16
$ws.Range("D16").Value = "60.667.18"
$ws.Range("E16").Value = "  +0.68%  "

# Row 17
$ws.Range("E17").Value = "  -2.28%  "

# Row 18
$ws.Range("D18").Value = "2.345.72"
$ws.Range("E18").Value = "  -1.75%  "

# Row 19
$ws.Range("D19").Value = "10.62"
$ws.Range("E19").Value = "  +0.66%  "

# Row 20
$ws.Range("D20").Value = "319.23"
$ws.Range("E20").Value = "  +2.07%  "

# Row 21
$ws.Range("D21").Value = "4.11"
$ws.Range("E21").Value = "  +1.68%  "

# Row 22
$ws.Range("D22").Value = "6.53"
$ws.Range("E22").Value = "  -2.15%  "

# Row 23
$ws.Range("D23").Value = "0.992"
$ws.Range("E23").Value = "  -0.76%  "

# Row 24
$ws.Range("D24").Value = "1.74"
$ws.Range("E24").Value = "  -2.34%  "

# Row 25
$ws.Range("D25").Value = "63.17"
$ws.Range("E25").Value = "  +0.96%  "

# Row 26
$ws.Range("D26").Value = "8.30"
$ws.Range("E26").Value = "  +9.13%  "

# Row 27
$ws.Range("D27").Value = "7.94"
$ws.Range("E27").Value = "  -0.35%  "

# Row 28
$ws.Range("D28").Value = "497.85"
$ws.Range("E28").Value = "  -0.42%  "

# Row 29
$ws.Range("E29").Value = "  -2.61%  "

# Row 30
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "0.0₃0866"
$ws.Range("E30").Value = "  -4.07%  "

# Row 31
$ws.Range("B31").Value = "Kaspa"
$ws.Range("C31").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D31").Value = "0.146"
$ws.Range("E31").Value = "  +2.06%  "

# Row 32
$ws.Range("D32").Value = "1.79"
$ws.Range("E32").Value = "  -1.72%  "

# Row 33
$ws.Range("D33").Value = "1.50"
$ws.Range("E33").Value = "  -3.26%  "

# Row 34
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  +0.03%  "

# Row 35
$ws.Range("D35").Value = "4.59"
$ws.Range("E35").Value = "  -0.47%  "

# Row 36
$ws.Range("D36").Value = "0.376"
$ws.Range("E36").Value = "  +1.34%  "

# Row 37
$ws.Range("E37").Value = "  +3.64%  "

# Row 38
$ws.Range("D38").Value = "5.27"
$ws.Range("E38").Value = "  -2.61%  "

# Row 39
$ws.Range("E39").Value = "  +7.27%  "

# Row 40
$ws.Range("D40").Value = "141.49"
$ws.Range("E40").Value = "  +3.16%  "

# Row 41
$ws.Range("D41").Value = "1.02"
$ws.Range("E41").Value = "  +2.16%  "

# Row 42
$ws.Range("D42").Value = "40.55"
$ws.Range("E42").Value = "  +0.62%  "

# Row 43
$ws.Range("D43").Value = "142.37"
$ws.Range("E43").Value = "  +2.05%  "

# Row 44
$ws.Range("D44").Value = "3.56"
$ws.Range("E44").Value = "  +1.58%  "

# Row 45
$ws.Range("D45").Value = "2.06"
$ws.Range("E45").Value = "  -4.34%  "

# Row 46
$ws.Range("D46").Value = "0.0515"
$ws.Range("E46").Value = "  +0.40%  "

# Row 47
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "0.570"
$ws.Range("E47").Value = "  -0.74%  "

# Row 48
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "19.05"
$ws.Range("E48").Value = "  -4.96%  "

# Row 49
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").Value = "0.0902"
$ws.Range("E49").Value = "  -2.03%  "

# Row 50
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "0.0220"
$ws.Range("E50").Value = "  -0.97%  "

# Row 51
$ws.Range("D51").Value = "16.24"
$ws.Range("E51").Value = "  -2.40%  "
